{"js": "// Elementalist card-sheet update.\n// Strategy: each logical change is applied as a literal find-and-replace\n// against the visible paragraph text (Word's search engine matches across\n// run/proofErr boundaries, so we can target human-readable phrases).\n\nasync function replaceOnce(body, searchText, replacement) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + searchText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacement, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1) All 16 \"summon command\" reminder-text cells: distance 4 -> distance 3\nawait replaceOnce(\n  body,\n  \"You may move two spaces OR command each of your summons within distance 4 to perform an ability.\",\n  \"You may move two spaces OR command each of your summons within distance 3 to perform an ability.\"\n);\n\n// ---- Familiars row ----\n\n// Fire Familiar: health 4 -> 6\nawait replaceOnce(\n  body,\n  \"Summon a Fire Familiar with 4 health on an\",\n  \"Summon a Fire Familiar with 6 health on an\"\n);\n\n// Fire Familiar ability 2: within 3 spaces -> within 2 spaces\nawait replaceOnce(\n  body,\n  \"2) Incinerate all enemies within 3 spaces for 2 damage.\",\n  \"2) Incinerate all enemies within 2 spaces for 2 damage.\"\n);\n\n// Ice Familiar: health 5 -> 6\nawait replaceOnce(\n  body,\n  \"Summon an Ice Familiar with 5 health on an\",\n  \"Summon an Ice Familiar with 6 health on an\"\n);\n\n// Ice Familiar ability 2: wording unchanged, but merge the proofErr-split runs\n// into one clean run (matches the source edit's run cleanup).\nawait replaceOnce(\n  body,\n  \"2) Freeze target enemy within 2 spaces for 2 damage; enemy cannot move this turn.\",\n  \"2) Freeze target enemy within 2 spaces for 2 damage; enemy cannot move this turn.\"\n);\n\n// Lightning Familiar: health 4 -> 6\nawait replaceOnce(\n  body,\n  \"Lightning Familiar with 4 health on an\",\n  \"Lightning Familiar with 6 health on an\"\n);\n\n// Lightning Familiar ability 2: within 3 spaces -> within 2 spaces, add trailing period\nawait replaceOnce(\n  body,\n  \"Shock an enemy within 3 spaces for 3 damage\",\n  \"Shock an enemy within 2 spaces for 3 damage.\"\n);\n\n// Earth Familiar: health 6 -> 8\nawait replaceOnce(\n  body,\n  \"Earth Familiar with 6 health on an\",\n  \"Earth Familiar with 8 health on an\"\n);\n\n// ---- Golems row ----\n\n// Fire Golem: health 8 -> 9\nawait replaceOnce(\n  body,\n  \"Summon a Fire Golem with health 8 on an adjacent space w\",\n  \"Summon a Fire Golem with health 9 on an adjacent space w\"\n);\n\n// Fire Golem ability 3: within 3 spaces -> within 2 spaces\nawait replaceOnce(\n  body,\n  \"Incinerate all enemies within 3 spaces for 3 damage\",\n  \"Incinerate all enemies within 2 spaces for 3 damage\"\n);\n\n// Ice Golem: health 10 -> 9\nawait replaceOnce(\n  body,\n  \"Ice Golem with health 10 on an adjacent space w\",\n  \"Ice Golem with health 9 on an adjacent space w\"\n);\n\n// Ice Golem ability 3: damage 3 -> 4\nawait replaceOnce(\n  body,\n  \"Freeze target enemy within 2 spaces for 3 damage; enemy cannot move this turn.\",\n  \"Freeze target enemy within 2 spaces for 4 damage; enemy cannot move this turn.\"\n);\n\n// Lightning Golem: health 8 -> 9\nawait replaceOnce(\n  body,\n  \"Lightning Golem with health 8 on an adjacent space w\",\n  \"Lightning Golem with health 9 on an adjacent space w\"\n);\n\n// Lightning Golem ability 2: within 3 spaces -> within 2 spaces, damage 5 -> 6\nawait replaceOnce(\n  body,\n  \"Shock an enemy within 3 spaces for 5 damage\",\n  \"Shock an enemy within 2 spaces for 6 damage\"\n);\n\n// ---- Merge-only cleanups (proofErr split runs -> plain text, no wording change) ----\n// \"Smash adjacent enemy: 5 damage\" occurs for Fire/Ice/Earth Golems.\nawait replaceOnce(body, \"ash adjacent enemy: 5 damage\", \"ash adjacent enemy: 5 damage\");\n\n// ---- Ongoing (\"Command\") cards ----\n\nawait replaceOnce(\n  body,\n  \"If a target is damaged by both 'Fire' and 'Lightning' this turn, each of their Ongoing effects disintegrate.\",\n  \"If a target is damaged by both 'Fire' and 'Lightning' this turn, strip off an Ongoing of your choice, and you may play a non-'Command' card.\"\n);\n\nawait replaceOnce(\n  body,\n  \"If a target is damaged by both 'Lightning' and 'Earth' this turn, they have a 50% chance to lose their next turn.\",\n  \"If a target is damaged by both 'Lightning' and 'Earth' this turn, You may play a non-'Command' card, and heal yourself 2 HP.\"\n);\n\nawait replaceOnce(\n  body,\n  \"If a target is damaged by both 'Fire' and 'Earth' this turn, the molten ground under their feet hardens, and they cannot move this turn or next.\",\n  \"If a target is damaged by both 'Fire' and 'Earth' this turn, strip off an Ongoing of your choice, and heal yourself 2 HP.\"\n);\n\n// \"Icy Hot\" -> \"Of Ice and Fire\" (title)\nawait replaceOnce(body, \"Icy Hot\", \"Of Ice and Fire\");\n\nawait replaceOnce(\n  body,\n  \"If a target is damaged by both 'Ice' and 'Fire' this turn, you heal for 4 health.\",\n  \"If a target is damaged by both 'Ice' and 'Fire' this turn, strip off an Ongoing of your choice, and all their cards next turn have Speed 10.\"\n);\n\nawait replaceOnce(\n  body,\n  \"If a target is damaged by both 'Ice' and 'Earth' this turn, all their cards next turn have Speed 10.\",\n  \"If a target is damaged by both 'Ice' and 'Earth' this turn, all their cards next turn have Speed 10, and heal yourself 2 HP.\"\n);\n\nawait replaceOnce(\n  body,\n  \"If a target is damaged by both 'Ice' and 'Lightning' this turn, you may play another non-'Command' card.\",\n  \"If a target is damaged by both 'Ice' and 'Lightning' this turn, all their cards next turn have Speed 10, and you may play a non-'Command' card.\"\n);\n", "ps1": "# Elementalist card-sheet update.\n# Strategy: each logical change is applied as a literal Find/Replace (wdReplaceAll)\n# against the document's content Range. Word's Find engine matches across\n# run/proofErr boundaries, so targeting human-readable phrases is reliable\n# even though the underlying text is split across several <w:r> runs.\n\nfunction Replace-Text($doc, $searchText, $replaceText) {\n    $find = $doc.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $replaceText\n    $result = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $result) {\n        throw \"Replace-Text: no match found for: $searchText\"\n    }\n}\n\n$d = $word.ActiveDocument\n\n# 1) All 16 \"summon command\" reminder-text cells: distance 4 -> distance 3\nReplace-Text $d \"You may move two spaces OR command each of your summons within distance 4 to perform an ability.\" \"You may move two spaces OR command each of your summons within distance 3 to perform an ability.\"\n\n# ---- Familiars row ----\n\n# Fire Familiar: health 4 -> 6\nReplace-Text $d \"Summon a Fire Familiar with 4 health on an\" \"Summon a Fire Familiar with 6 health on an\"\n\n# Fire Familiar ability 2: within 3 spaces -> within 2 spaces\nReplace-Text $d \"2) Incinerate all enemies within 3 spaces for 2 damage.\" \"2) Incinerate all enemies within 2 spaces for 2 damage.\"\n\n# Ice Familiar: health 5 -> 6\nReplace-Text $d \"Summon an Ice Familiar with 5 health on an\" \"Summon an Ice Familiar with 6 health on an\"\n\n# Ice Familiar ability 2: wording unchanged, but merge the proofErr-split runs\n# into one clean run (matches the source edit's run cleanup).\nReplace-Text $d \"2) Freeze target enemy within 2 spaces for 2 damage; enemy cannot move this turn.\" \"2) Freeze target enemy within 2 spaces for 2 damage; enemy cannot move this turn.\"\n\n# Lightning Familiar: health 4 -> 6\nReplace-Text $d \"Lightning Familiar with 4 health on an\" \"Lightning Familiar with 6 health on an\"\n\n# Lightning Familiar ability 2: within 3 spaces -> within 2 spaces, add trailing period\nReplace-Text $d \"Shock an enemy within 3 spaces for 3 damage\" \"Shock an enemy within 2 spaces for 3 damage.\"\n\n# Earth Familiar: health 6 -> 8\nReplace-Text $d \"Earth Familiar with 6 health on an\" \"Earth Familiar with 8 health on an\"\n\n# ---- Golems row ----\n\n# Fire Golem: health 8 -> 9\nReplace-Text $d \"Summon a Fire Golem with health 8 on an adjacent space w\" \"Summon a Fire Golem with health 9 on an adjacent space w\"\n\n# Fire Golem ability 3: within 3 spaces -> within 2 spaces\nReplace-Text $d \"Incinerate all enemies within 3 spaces for 3 damage\" \"Incinerate all enemies within 2 spaces for 3 damage\"\n\n# Ice Golem: health 10 -> 9\nReplace-Text $d \"Ice Golem with health 10 on an adjacent space w\" \"Ice Golem with health 9 on an adjacent space w\"\n\n# Ice Golem ability 3: damage 3 -> 4\nReplace-Text $d \"Freeze target enemy within 2 spaces for 3 damage; enemy cannot move this turn.\" \"Freeze target enemy within 2 spaces for 4 damage; enemy cannot move this turn.\"\n\n# Lightning Golem: health 8 -> 9\nReplace-Text $d \"Lightning Golem with health 8 on an adjacent space w\" \"Lightning Golem with health 9 on an adjacent space w\"\n\n# Lightning Golem ability 2: within 3 spaces -> within 2 spaces, damage 5 -> 6\nReplace-Text $d \"Shock an enemy within 3 spaces for 5 damage\" \"Shock an enemy within 2 spaces for 6 damage\"\n\n# ---- Merge-only cleanups (proofErr split runs -> plain text, no wording change) ----\n# \"Smash adjacent enemy: 5 damage\" occurs for Fire/Ice/Earth Golems.\nReplace-Text $d \"ash adjacent enemy: 5 damage\" \"ash adjacent enemy: 5 damage\"\n\n# ---- Ongoing (\"Command\") cards ----\n\nReplace-Text $d \"If a target is damaged by both 'Fire' and 'Lightning' this turn, each of their Ongoing effects disintegrate.\" \"If a target is damaged by both 'Fire' and 'Lightning' this turn, strip off an Ongoing of your choice, and you may play a non-'Command' card.\"\n\nReplace-Text $d \"If a target is damaged by both 'Lightning' and 'Earth' this turn, they have a 50% chance to lose their next turn.\" \"If a target is damaged by both 'Lightning' and 'Earth' this turn, You may play a non-'Command' card, and heal yourself 2 HP.\"\n\nReplace-Text $d \"If a target is damaged by both 'Fire' and 'Earth' this turn, the molten ground under their feet hardens, and they cannot move this turn or next.\" \"If a target is damaged by both 'Fire' and 'Earth' this turn, strip off an Ongoing of your choice, and heal yourself 2 HP.\"\n\n# \"Icy Hot\" -> \"Of Ice and Fire\" (title)\nReplace-Text $d \"Icy Hot\" \"Of Ice and Fire\"\n\nReplace-Text $d \"If a target is damaged by both 'Ice' and 'Fire' this turn, you heal for 4 health.\" \"If a target is damaged by both 'Ice' and 'Fire' this turn, strip off an Ongoing of your choice, and all their cards next turn have Speed 10.\"\n\nReplace-Text $d \"If a target is damaged by both 'Ice' and 'Earth' this turn, all their cards next turn have Speed 10.\" \"If a target is damaged by both 'Ice' and 'Earth' this turn, all their cards next turn have Speed 10, and heal yourself 2 HP.\"\n\nReplace-Text $d \"If a target is damaged by both 'Ice' and 'Lightning' this turn, you may play another non-'Command' card.\" \"If a target is damaged by both 'Ice' and 'Lightning' this turn, all their cards next turn have Speed 10, and you may play a non-'Command' card.\"\n"}
